$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the "Tuna" Melon block (rows 49-52), pushing the
# existing data (previously rows 49-70) down to rows 53-74. This matches the
# resulting worksheet exactly as described in the diff (dimension A1:R70 -> A1:R74,
# with 4 brand-new "Calameño" price records inserted ahead of the "Tuna" records).
$ws.Range("A49:R52").Insert()

# New row 49: Calameño / Primera
$ws.Cells.Item(49, 1).Value = 8
$ws.Cells.Item(49, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(49, 5).Value = 4
$ws.Cells.Item(49, 6).Value = 100112027
$ws.Cells.Item(49, 7).Value = "Melón"
$ws.Cells.Item(49, 8).Value = "Calameño"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 5000
$ws.Cells.Item(49, 11).Value = 850
$ws.Cells.Item(49, 12).Value = 900
$ws.Cells.Item(49, 13).Value = 875
$ws.Cells.Item(49, 14).Value = "`$/unidad"
$ws.Cells.Item(49, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(49, 16).Value = 875
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# New row 50: Calameño / Extra
$ws.Cells.Item(50, 1).Value = 8
$ws.Cells.Item(50, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value = "Coquimbo"
$ws.Cells.Item(50, 4).Value = Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = 100112027
$ws.Cells.Item(50, 7).Value = "Melón"
$ws.Cells.Item(50, 8).Value = "Calameño"
$ws.Cells.Item(50, 9).Value = "Extra"
$ws.Cells.Item(50, 10).Value = 4400
$ws.Cells.Item(50, 11).Value = 1000
$ws.Cells.Item(50, 12).Value = 1100
$ws.Cells.Item(50, 13).Value = 1050
$ws.Cells.Item(50, 14).Value = "`$/unidad"
$ws.Cells.Item(50, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 16).Value = 1050
$ws.Cells.Item(50, 17).Value = 1
$ws.Cells.Item(50, 18).Value = "Hortaliza"

# New row 51: Calameño / Primera
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = 100112027
$ws.Cells.Item(51, 7).Value = "Melón"
$ws.Cells.Item(51, 8).Value = "Calameño"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 4000
$ws.Cells.Item(51, 11).Value = 800
$ws.Cells.Item(51, 12).Value = 900
$ws.Cells.Item(51, 13).Value = 850
$ws.Cells.Item(51, 14).Value = "`$/unidad"
$ws.Cells.Item(51, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 16).Value = 850
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# New row 52: Calameño / Super
$ws.Cells.Item(52, 1).Value = 8
$ws.Cells.Item(52, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = 100112027
$ws.Cells.Item(52, 7).Value = "Melón"
$ws.Cells.Item(52, 8).Value = "Calameño"
$ws.Cells.Item(52, 9).Value = "Super"
$ws.Cells.Item(52, 10).Value = 4800
$ws.Cells.Item(52, 11).Value = 1200
$ws.Cells.Item(52, 12).Value = 1300
$ws.Cells.Item(52, 13).Value = 1250
$ws.Cells.Item(52, 14).Value = "`$/unidad"
$ws.Cells.Item(52, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 16).Value = 1250
$ws.Cells.Item(52, 17).Value = 1
$ws.Cells.Item(52, 18).Value = "Hortaliza"
